$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "ABOUT ME" paragraph - reword the closing sentence and move
# the "_GoBack" bookmark (Word's "last edit" marker) into this paragraph,
# right after the newly typed "to apply".
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "I am seeking to be part of SaasyCloud and contribute towards achieving",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I am seeking to apply my skills and abilities towards achieving", 2)

# Re-locate the new pieces of text so we can split them back into
# individually-typed runs (mirrors how Word leaves separate runs behind
# after a sequence of edits) and drop the bookmark at the right spot.
$rToApply = $d.Content.Duplicate
$rToApply.Find.Execute("to apply", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryAfterToApply = $rToApply.End

$rTowards = $d.Content.Duplicate
$rTowards.Find.Execute("towards achieving", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryBeforeTowards = $rTowards.Start

$boundaryAfterSpace = $boundaryAfterToApply + 1

# Split the run into separate pieces using temporary bookmarks, then
# drop the temporary ones, leaving a real "_GoBack" bookmark exactly
# where Word would have left it after typing "to apply".
$d.Bookmarks.Add("zzTmpSplit2", $d.Range($boundaryBeforeTowards, $boundaryBeforeTowards))
$d.Bookmarks.Add("zzTmpSplit1", $d.Range($boundaryAfterSpace, $boundaryAfterSpace))
$d.Bookmarks.Add("_GoBack", $d.Range($boundaryAfterToApply, $boundaryAfterToApply))
$d.Bookmarks("zzTmpSplit2").Delete()
$d.Bookmarks("zzTmpSplit1").Delete()

# -----------------------------------------------------------------------
# Change 2: "EDUCATION" heading no longer starts a rendered page, so the
# cached w:lastRenderedPageBreak marker on its run should disappear.
# Re-applying the same text through Find/Replace rebuilds the run and
# drops the stale rendering marker.
# -----------------------------------------------------------------------
$d.Content.Find.Execute("EDUCATION", $true, $false, $false, $false, $false, $true, 1, $false, "EDUCATION", 2)

# -----------------------------------------------------------------------
# Change 3: "Worked with the lead..." bullet - the bookmark that used to
# sit mid-sentence is gone and the sentence is now one contiguous run.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Worked with the lead to adjust how the tickets were handled",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Worked with the lead to adjust how the tickets were handled", 2)
